# Auto-generated: apply scheduled-runner market-price updates to Sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1398.875
$ws.Range("I62").Value = 1300.5294
$ws.Range("J62").Value = 1637.7142
$ws.Range("K62").Value = 1300.5294
$ws.Range("L62").Value = 1637.7142
$ws.Range("M62").Value = -676.5293999999999
$ws.Range("N62").Value = -2885.7142
$ws.Range("H65").Value = 1398.875
$ws.Range("I65").Value = 1300.5294
$ws.Range("J65").Value = 1637.7142
$ws.Range("K65").Value = 6502.646999999999
$ws.Range("L65").Value = 8188.571
$ws.Range("M65").Value = -3382.646999999999
$ws.Range("N65").Value = -14428.571
$ws.Range("H113").Value = 3923865.8
$ws.Range("I113").Value = 4764087
$ws.Range("J113").Value = 2833.3333
$ws.Range("K113").Value = 4764087
$ws.Range("L113").Value = 2833.3333
$ws.Range("M113").Value = -4760833
$ws.Range("N113").Value = -9341.3333
$ws.Range("H116").Value = 5758.9287
$ws.Range("I116").Value = 6613
$ws.Range("K116").Value = 6613
$ws.Range("M116").Value = -3171
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H132").Value = 1525.5294
$ws.Range("I132").Value = 1616.7858
$ws.Range("J132").Value = 1099.6666
$ws.Range("K132").Value = 4850.357400000001
$ws.Range("L132").Value = 3298.9998
$ws.Range("M132").Value = -2320.357400000001
$ws.Range("N132").Value = -8358.9998
$ws.Range("H141").Value = 1930
$ws.Range("I141").Value = 1447.5
$ws.Range("J141").Value = 2787.7778
$ws.Range("K141").Value = 4342.5
$ws.Range("L141").Value = 8363.3334
$ws.Range("M141").Value = 837.5
$ws.Range("N141").Value = -18723.3334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6621.1953
$ws.Range("I32").Value = 3060.543
$ws.Range("J32").Value = 27391.666
$ws.Range("K32").Value = 3060.543
$ws.Range("L32").Value = 27391.666
$ws.Range("M32").Value = -2773.543
$ws.Range("N32").Value = -27965.666
$ws.Range("H61").Value = 2777.158
$ws.Range("I61").Value = 2115.375
$ws.Range("J61").Value = 6306.6665
$ws.Range("K61").Value = 2115.375
$ws.Range("L61").Value = 6306.6665
$ws.Range("M61").Value = -1903.375
$ws.Range("N61").Value = -6730.6665
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H136").Value = 2777.158
$ws.Range("I136").Value = 2115.375
$ws.Range("J136").Value = 6306.6665
$ws.Range("K136").Value = 6346.125
$ws.Range("L136").Value = 18919.9995
$ws.Range("M136").Value = -3796.125
$ws.Range("N136").Value = -24019.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1504.3529
$ws.Range("I86").Value = 1459.8096
$ws.Range("J86").Value = 1712.2222
$ws.Range("K86").Value = 1459.8096
$ws.Range("L86").Value = 1712.2222
$ws.Range("M86").Value = -336.8096
$ws.Range("N86").Value = -3958.2222
$ws.Range("H89").Value = 1504.3529
$ws.Range("I89").Value = 1459.8096
$ws.Range("J89").Value = 1712.2222
$ws.Range("K89").Value = 7299.048000000001
$ws.Range("L89").Value = 8561.110999999999
$ws.Range("M89").Value = -1683.048000000001
$ws.Range("N89").Value = -19793.111
$ws.Range("H99").Value = 47619984
$ws.Range("I99").Value = 58824388
$ws.Range("J99").Value = 1274.75
$ws.Range("K99").Value = 58824388
$ws.Range("L99").Value = 1274.75
$ws.Range("M99").Value = -58822890
$ws.Range("N99").Value = -4270.75
$ws.Range("H100").Value = 30643
$ws.Range("J100").Value = 30643
$ws.Range("L100").Value = 30643
$ws.Range("N100").Value = -32807
$ws.Range("H105").Value = 9010.833000000001
$ws.Range("I105").Value = 12296.158
$ws.Range("J105").Value = 3336.182
$ws.Range("K105").Value = 12296.158
$ws.Range("L105").Value = 3336.182
$ws.Range("M105").Value = -10549.158
$ws.Range("N105").Value = -6830.182

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1792.3334
$ws.Range("I16").Value = 1487.9231
$ws.Range("K16").Value = 1487.9231
$ws.Range("M16").Value = -1200.9231
$ws.Range("H31").Value = 3884.5151
$ws.Range("I31").Value = 1521.9565
$ws.Range("J31").Value = 9318.4
$ws.Range("K31").Value = 1521.9565
$ws.Range("L31").Value = 9318.4
$ws.Range("M31").Value = -1226.9565
$ws.Range("N31").Value = -9908.4
$ws.Range("H34").Value = 3884.5151
$ws.Range("I34").Value = 1521.9565
$ws.Range("J34").Value = 9318.4
$ws.Range("K34").Value = 1521.9565
$ws.Range("L34").Value = 9318.4
$ws.Range("M34").Value = -1319.9565
$ws.Range("N34").Value = -9722.4
$ws.Range("H94").Value = 4541.7827
$ws.Range("J94").Value = 4357
$ws.Range("L94").Value = 4357
$ws.Range("N94").Value = -5259
$ws.Range("H96").Value = 14719.8
$ws.Range("J96").Value = 14719.8
$ws.Range("L96").Value = 14719.8
$ws.Range("N96").Value = -20211.8
$ws.Range("H99").Value = 29284.6
$ws.Range("I99").Value = 29102.25
$ws.Range("J99").Value = 30014
$ws.Range("K99").Value = 29102.25
$ws.Range("L99").Value = 30014
$ws.Range("M99").Value = -27604.25
$ws.Range("N99").Value = -33010
$ws.Range("H105").Value = 1712.55
$ws.Range("I105").Value = 1675.1111
$ws.Range("J105").Value = 2049.5
$ws.Range("K105").Value = 1675.1111
$ws.Range("L105").Value = 2049.5
$ws.Range("M105").Value = 71.88889999999992
$ws.Range("N105").Value = -5543.5
$ws.Range("H113").Value = 1792.3334
$ws.Range("I113").Value = 1487.9231
$ws.Range("K113").Value = 1487.9231
$ws.Range("M113").Value = 682.0769
$ws.Range("H126").Value = 29284.6
$ws.Range("I126").Value = 29102.25
$ws.Range("J126").Value = 30014
$ws.Range("K126").Value = 87306.75
$ws.Range("L126").Value = 90042
$ws.Range("M126").Value = -84836.75
$ws.Range("N126").Value = -94982
$ws.Range("H134").Value = 2525.0466
$ws.Range("I134").Value = 2466.5898
$ws.Range("K134").Value = 7399.769400000001
$ws.Range("M134").Value = -4864.769400000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4595176.5
$ws.Range("I122").Value = 4053179.5
$ws.Range("J122").Value = 5558726.5
$ws.Range("K122").Value = 12159538.5
$ws.Range("L122").Value = 16676179.5
$ws.Range("M122").Value = -12157088.5
$ws.Range("N122").Value = -16681079.5
$ws.Range("H123").Value = 23627.572
$ws.Range("J123").Value = 23627.572
$ws.Range("L123").Value = 23627.572
$ws.Range("N123").Value = -28527.572
$ws.Range("H126").Value = 6240.375
$ws.Range("I126").Value = 8436.933999999999
$ws.Range("K126").Value = 25310.802
$ws.Range("M126").Value = -22840.802

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 14320
$ws.Range("J104").Value = 14320
$ws.Range("L104").Value = 14320
$ws.Range("N104").Value = -21308
$ws.Range("H122").Value = 2397023.8
$ws.Range("I122").Value = 3573229.5
$ws.Range("J122").Value = 716729.6
$ws.Range("K122").Value = 10719688.5
$ws.Range("L122").Value = 2150188.8
$ws.Range("M122").Value = -10717238.5
$ws.Range("N122").Value = -2155088.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1290.2368
$ws.Range("I132").Value = 994.6896400000001
$ws.Range("K132").Value = 2984.06892
$ws.Range("M132").Value = -454.0689200000002
